# Changing Text Case (UPPER, LOWER, PROPER)
# - Column D ("Full Name"): wrap the existing CONCAT(...) formula in PROPER()
#   so names always render in proper case.
# - Column E ("Email"): wrap the existing "C&.&B&@pushpin.com" formula in
#   LOWER() so email addresses always render in lower case.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Top (non-shared) formulas on row 4
$ws.Range("D4").Formula = '=PROPER(CONCAT(C4," ",B4))'
$ws.Range("E4").Formula = '=LOWER(C4&"."&B4&"@pushpin.com")'

# Shared formula block rows 5:38 - setting the range fills every cell,
# and Excel will re-establish the shared formula group starting at D5/E5.
$ws.Range("D5:D38").Formula = '=PROPER(CONCAT(C5," ",B5))'
$ws.Range("E5:E38").Formula = '=LOWER(C5&"."&B5&"@pushpin.com")'

# Reflect the author's final selection (E4:E38) as left by the edit.
$ws.Range("E4:E38").Select()
